$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "thaovy"
$ws.Range("B3").Value = 123

$ws.Range("D4").Select()
